$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: filler/separator row (copy border+font formatting from row 19 pattern: s=6/7)
$ws.Range("A19:E19").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("A22").Value = 'SCRIPT/G01P03A/um2504.ssb'

# Row 23: data row (copy formatting from row 21 pattern: s=4/5)
$ws.Range("A21:E21").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)
$ws.Range("C23").Value = ' Eep! Y-y-yeesh…'
$ws.Range("A23").Value = 'SCRIPT/T01P02A/us0111.ssb'
$ws.Range("B23").Value = 206
$ws.Range("D23").Value = ' Оой! Д-д-дааа...'
$ws.Range("E23").Value = ' Ïïê! Ä-ä-äààà...'

# Row 24: filler/separator row (copy formatting from row 19 pattern: s=6/7)
$ws.Range("A19:E19").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)
$ws.Range("A24").Value = 'SCRIPT/T01P02A/us3103.ssb'

# Row 25: data row (copy formatting from row 21 pattern: s=4/5)
$ws.Range("A21:E21").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("C25").Value = ' Huh? Why is it so empty\nin here?'
$ws.Range("A25").Value = 'SCRIPT/P01P04A/us3120.ssb'
$ws.Range("B25").Value = 225
$ws.Range("D25").Value = ' Что? Почему здесь почти никого\nнет?'
$ws.Range("E25").Value = ' Œóï? Ðïœåíô èäåòû ðïœóé îéëïãï\nîåó?'

# Row heights (44.2pt / 43.2pt wrap rows, matching rows above)
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 43.2
$ws.Rows.Item(24).RowHeight = 43.2
$ws.Rows.Item(25).RowHeight = 43.2

$excel.Application.GoTo($ws.Range("A25"))
$ws.Range("D27").Select()

Write-Output "done"
